$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1681445254291393
$ws.Range("C2").Value = 0.04385043669510873
$ws.Range("B3").Value = 0.3452536934207687
$ws.Range("C3").Value = -0.1566425969469884
$ws.Range("B4").Value = 0.3556355745522564
$ws.Range("C4").Value = 0.1990384573524147
$ws.Range("B5").Value = 0.431523411195365
$ws.Range("C5").Value = -0.3941062953087233
$ws.Range("B6").Value = 0.360484414827775
$ws.Range("C6").Value = -0.02804511262863444
$ws.Range("B7").Value = 0.3972390355134647
$ws.Range("C7").Value = -0.007063301683255129
$ws.Range("B8").Value = 0.2321230097122929
$ws.Range("C8").Value = 0.5083813793829081
$ws.Range("B9").Value = 0.1872709057121062
$ws.Range("C9").Value = 0.3528597334377797
$ws.Range("B10").Value = 0.2154501347641879
$ws.Range("C10").Value = -0.4616186901550188
$ws.Range("B11").Value = 0.1398838869808469
$ws.Range("C11").Value = 0.2414382423752041
$ws.Range("B12").Value = -0.04096520384890977
$ws.Range("C12").Value = -0.1502088144699407
$ws.Range("B13").Value = 0.1113490370755945
$ws.Range("C13").Value = -0.1770168458584153
$ws.Range("B14").Value = -0.04382704684561507
$ws.Range("C14").Value = -0.1946009543459019
$ws.Range("B15").Value = 0.1445655189723084
$ws.Range("C15").Value = 0.01494890086795084
$ws.Range("B16").Value = 0.2439679294167482
$ws.Range("C16").Value = 0.1772231259449329
$ws.Range("B17").Value = -0.02711960050508749
$ws.Range("C17").Value = -0.004815648980808656
